$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.103994727134705
$ws.Range("B1").Value = 1.778860211372375
$ws.Range("C1").Value = 9.218070983886719
$ws.Range("D1").Value = 2.399453163146973
$ws.Range("E1").Value = 1.288108587265015
